$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "10. Create Company page" checklist item. The items that
# used to be numbered 10-12 are renumbered to 11-13, and the whole list
# keeps its row positions (rows 13-16), with row 16 being the brand new
# row for the former item 12 ("Mentor--Student linking Page"), now #13.
$ws.Range("A13").Value = "10. Create Company page"
$ws.Range("A14").Value = "11. Account creation page"
$ws.Range("A15").Value = "12. Change Account info page"
$ws.Range("A16").Value = "13. Mentor--Student linking Page"

# Mark the new "Create Company page" / "Account creation page" rows, as
# well as "3. Computerize revised ERD", as accomplished.
$ws.Range("B13").Value = "X"
$ws.Range("B14").Value = "X"
$ws.Range("B6").Value = "X"

# Center-align the whole "Accomplished?" column for every checklist row,
# including the blank spacer row (3) and the new trailing blank row (17).
$ws.Range("B3:B17").HorizontalAlignment = -4108

# Move the active selection to B10, matching the saved view state.
$ws.Range("B10").Select()
